$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds date-like text ("2025-11-14"); force Text format first so
# Excel does not silently coerce the literal into a date serial number.
$ws.Range("B3:B6").NumberFormat = "@"

# Row 3
$ws.Range("A3").Value = "Danish 2nd Division"
$ws.Range("B3").Value = "2025-11-14"
$ws.Range("C3").Value = "14:00:00"
$ws.Range("D3").Value = "HIK Hellerup"
$ws.Range("E3").Value = "Vendsyssel FF"
$ws.Range("F3").Value = 3.1
$ws.Range("G3").Value = 3.95
$ws.Range("H3").Value = 2.24
$ws.Range("I3").Value = 2.58
$ws.Range("J3").Value = 3.2
$ws.Range("K3").Value = 3.8
$ws.Range("L3").Value = 1.34
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 3.35
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 1.83
$ws.Range("Q3").Value = 1.98
$ws.Range("R3").Value = 1.31
$ws.Range("S3").Value = 3.5
$ws.Range("T3").Value = 1.75
$ws.Range("U3").Value = 2.04
$ws.Range("V3").Value = 1.63
$ws.Range("W3").Value = 1.36
$ws.Range("X3").Value = 16.5
$ws.Range("Y3").Value = 12
$ws.Range("Z3").Value = 19
$ws.Range("AA3").Value = 40
$ws.Range("AB3").Value = 14.5
$ws.Range("AC3").Value = 9.4
$ws.Range("AD3").Value = 13.5
$ws.Range("AE3").Value = 34
$ws.Range("AF3").Value = 28
$ws.Range("AG3").Value = 17
$ws.Range("AH3").Value = 22
$ws.Range("AI3").Value = 50
$ws.Range("AJ3").Value = 75
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 65
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 50
$ws.Range("AO3").Value = 27

# Row 4
$ws.Range("A4").Value = "Spanish Segunda Division"
$ws.Range("B4").Value = "2025-11-14"
$ws.Range("C4").Value = "16:30:00"
$ws.Range("D4").Value = "Valladolid"
$ws.Range("E4").Value = "Las Palmas"
$ws.Range("F4").Value = 2.3
$ws.Range("G4").Value = 2.48
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 3.95
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 3.2
$ws.Range("L4").Value = 1.52
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 2.66
$ws.Range("O4").Value = 1.51
$ws.Range("P4").Value = 1.57
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.2
$ws.Range("S4").Value = 5.1
$ws.Range("T4").Value = 2.02
$ws.Range("U4").Value = 1.79
$ws.Range("V4").Value = 1.33
$ws.Range("W4").Value = 1.68
$ws.Range("X4").Value = 9.6
$ws.Range("Y4").Value = 11
$ws.Range("Z4").Value = 25
$ws.Range("AA4").Value = 90
$ws.Range("AB4").Value = 8.4
$ws.Range("AC4").Value = 7.2
$ws.Range("AD4").Value = 17
$ws.Range("AE4").Value = 60
$ws.Range("AF4").Value = 14
$ws.Range("AG4").Value = 12.5
$ws.Range("AH4").Value = 25
$ws.Range("AI4").Value = 85
$ws.Range("AJ4").Value = 36
$ws.Range("AK4").Value = 36
$ws.Range("AL4").Value = 65
$ws.Range("AM4").Value = 200
$ws.Range("AN4").Value = 38
$ws.Range("AO4").Value = 85

# Row 5
$ws.Range("A5").Value = "Welsh Premiership"
$ws.Range("B5").Value = "2025-11-14"
$ws.Range("C5").Value = "16:45:00"
$ws.Range("D5").Value = "Flint Town United"
$ws.Range("E5").Value = "The New Saints"
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 16.5
$ws.Range("H5").Value = 1.22
$ws.Range("I5").Value = 1.31
$ws.Range("J5").Value = 6.2
$ws.Range("K5").Value = 8.4
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 7.8
$ws.Range("O5").Value = 1.11
$ws.Range("P5").Value = 3.3
$ws.Range("Q5").Value = 1.29
$ws.Range("R5").Value = 1.95
$ws.Range("S5").Value = 1.84
$ws.Range("T5").Value = 1.78
$ws.Range("U5").Value = 2.04
$ws.Range("V5").Value = 4.2
$ws.Range("W5").Value = 1.06
$ws.Range("X5").Value = 55
$ws.Range("Y5").Value = 18.5
$ws.Range("Z5").Value = 13.5
$ws.Range("AA5").Value = 13.5
$ws.Range("AB5").Value = 70
$ws.Range("AC5").Value = 22
$ws.Range("AD5").Value = 15
$ws.Range("AE5").Value = 16.5
$ws.Range("AF5").Value = 160
$ws.Range("AG5").Value = 55
$ws.Range("AH5").Value = 34
$ws.Range("AI5").Value = 36
$ws.Range("AJ5").Value = 480
$ws.Range("AK5").Value = 190
$ws.Range("AL5").Value = 130
$ws.Range("AM5").Value = 130
$ws.Range("AN5").Value = 170
$ws.Range("AO5").Value = 3.35

# Row 6
$ws.Range("A6").Value = "Welsh Premiership"
$ws.Range("B6").Value = "2025-11-14"
$ws.Range("C6").Value = "16:45:00"
$ws.Range("D6").Value = "Cardiff Metropolitan"
$ws.Range("E6").Value = "Briton Ferry Llansawel"
$ws.Range("F6").Value = 1.67
$ws.Range("G6").Value = 1.81
$ws.Range("H6").Value = 4.7
$ws.Range("I6").Value = 5.7
$ws.Range("J6").Value = 3.95
$ws.Range("K6").Value = 4.7
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 4.3
$ws.Range("O6").Value = 1.23
$ws.Range("P6").Value = 2.14
$ws.Range("Q6").Value = 1.71
$ws.Range("R6").Value = 1.44
$ws.Range("S6").Value = 2.8
$ws.Range("T6").Value = 1.72
$ws.Range("U6").Value = 2.12
$ws.Range("V6").Value = 1.21
$ws.Range("W6").Value = 2.22
$ws.Range("X6").Value = 19.5
$ws.Range("Y6").Value = 22
$ws.Range("Z6").Value = 44
$ws.Range("AA6").Value = 140
$ws.Range("AB6").Value = 10.5
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 21
$ws.Range("AE6").Value = 70
$ws.Range("AF6").Value = 12
$ws.Range("AG6").Value = 10
$ws.Range("AH6").Value = 19.5
$ws.Range("AI6").Value = 70
$ws.Range("AJ6").Value = 19
$ws.Range("AK6").Value = 18
$ws.Range("AL6").Value = 32
$ws.Range("AM6").Value = 110
$ws.Range("AN6").Value = 11
$ws.Range("AO6").Value = 70

# Restore default (Normal) style on column B so the saved XML has no
# stray style index on these cells (matches plain inlineStr cells elsewhere).
$ws.Range("B3:B6").Style = "Normal"
